$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at position 82, shifting existing rows 82-117 down to 86-121
$ws.Rows("82:85").Insert()

# Row 82
$ws.Range("A82").Value = 2
$ws.Range("B82").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C82").Value = "Coquimbo"
$ws.Range("D82").Value = 44553
$ws.Range("E82").Value = 4
$ws.Range("F82").Value = "Fruta"
$ws.Range("G82").Value = 100103
$ws.Range("H82").Value = "Frutos de hueso (carozo)"
$ws.Range("I82").Value = 100103001
$ws.Range("J82").Value = "Cereza"
$ws.Range("K82").Value = "Lapins"
$ws.Range("L82").Value = "Primera"
$ws.Range("M82").Value = 400
$ws.Range("N82").Value = 9500
$ws.Range("O82").Value = 10000
$ws.Range("P82").Value = 9750
$ws.Range("Q82").Value = "`$/bandeja 10 kilos"
$ws.Range("R82").Value = "Región de O'Higgins"
$ws.Range("S82").Value = 975
$ws.Range("T82").Value = 10

# Row 83
$ws.Range("A83").Value = 2
$ws.Range("B83").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C83").Value = "Coquimbo"
$ws.Range("D83").Value = 44553
$ws.Range("E83").Value = 4
$ws.Range("F83").Value = "Fruta"
$ws.Range("G83").Value = 100103
$ws.Range("H83").Value = "Frutos de hueso (carozo)"
$ws.Range("I83").Value = 100103001
$ws.Range("J83").Value = "Cereza"
$ws.Range("K83").Value = "Lapins"
$ws.Range("L83").Value = "Segunda"
$ws.Range("M83").Value = 320
$ws.Range("N83").Value = 7500
$ws.Range("O83").Value = 8000
$ws.Range("P83").Value = 7750
$ws.Range("Q83").Value = "`$/bandeja 10 kilos"
$ws.Range("R83").Value = "Región de O'Higgins"
$ws.Range("S83").Value = 775
$ws.Range("T83").Value = 10

# Row 84
$ws.Range("A84").Value = 2
$ws.Range("B84").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C84").Value = "Coquimbo"
$ws.Range("D84").Value = 44553
$ws.Range("E84").Value = 4
$ws.Range("F84").Value = "Fruta"
$ws.Range("G84").Value = 100103
$ws.Range("H84").Value = "Frutos de hueso (carozo)"
$ws.Range("I84").Value = 100103001
$ws.Range("J84").Value = "Cereza"
$ws.Range("K84").Value = "Rainier"
$ws.Range("L84").Value = "Primera"
$ws.Range("M84").Value = 440
$ws.Range("N84").Value = 17000
$ws.Range("O84").Value = 18000
$ws.Range("P84").Value = 17500
$ws.Range("Q84").Value = "`$/caja 15 kilos"
$ws.Range("R84").Value = "Provincia de Curicó"
$ws.Range("S84").Value = 1167
$ws.Range("T84").Value = 15

# Row 85
$ws.Range("A85").Value = 2
$ws.Range("B85").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C85").Value = "Coquimbo"
$ws.Range("D85").Value = 44553
$ws.Range("E85").Value = 4
$ws.Range("F85").Value = "Fruta"
$ws.Range("G85").Value = 100103
$ws.Range("H85").Value = "Frutos de hueso (carozo)"
$ws.Range("I85").Value = 100103001
$ws.Range("J85").Value = "Cereza"
$ws.Range("K85").Value = "Santina"
$ws.Range("L85").Value = "Primera"
$ws.Range("M85").Value = 400
$ws.Range("N85").Value = 9500
$ws.Range("O85").Value = 10000
$ws.Range("P85").Value = 9750
$ws.Range("Q85").Value = "`$/bandeja 10 kilos"
$ws.Range("R85").Value = "Región de O'Higgins"
$ws.Range("S85").Value = 975
$ws.Range("T85").Value = 10
